$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-11 09:49:21"
$wsZh.Range("H2").Value = "2016-03-11 09:49:37"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-11 09:49:24"
$wsDe.Range("H2").Value = "2016-03-11 09:49:42"
